# Attendance report sync: reverse the order of the comma-separated
# "Recorded By" entries in column G for every data row, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com" ->  "backup@backdoor.com, System"
# A handful of rows were already in the target order upstream and were
# left untouched by the source sync, so we skip those explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value must stay exactly as-is (not reversed).
$exceptionRows = @(7, 34, 61)

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    if ($exceptionRows -contains $r) {
        continue
    }

    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -eq $null) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    # Build the reversed list (avoid [array]::Reverse - it does not
    # mutate in place in this runtime).
    $n = $parts.Count
    $reversedParts = @()
    for ($i = $n - 1; $i -ge 0; $i--) {
        $reversedParts += $parts[$i]
    }

    $newVal = [string]::Join(", ", $reversedParts)
    $cell.Value = $newVal
}
